$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in the header (F1): 13:00 -> 13:15
$ws.Range("F1").Value = "Last status check on: 25.02.2022 13:15"

# Row 7 (MOL Olomoucka) was updated by the AWS bash cron job with fresh
# numeric data instead of the previous placeholder text values.
$ws.Range("D7").Value = 0.21

$ws.Range("E7").Value = 44617.54188657407
$ws.Range("E7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
